$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    # Force text format so numeric-looking strings (e.g. "22.80", "0.900")
    # keep their exact original formatting instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value2 = $Text
    # Restore the default style so no stray formatting is introduced.
    $cell.Style = "Normal"
}

Set-TextValue "D2" "27.440.39"
Set-TextValue "E2" "  +0.25%  "
Set-TextValue "D3" "1.633.34"
Set-TextValue "E3" "  -1.02%  "
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "212.22"
Set-TextValue "E5" "  -0.52%  "
Set-TextValue "E6" "  +4.61%  "
Set-TextValue "E7" "  +0.11%  "
Set-TextValue "D8" "22.80"
Set-TextValue "E8" "  -5.69%  "
Set-TextValue "E9" "  -2.48%  "
Set-TextValue "E10" "  -0.83%  "
Set-TextValue "E11" "  +1.23%  "
Set-TextValue "D12" "1.866.03"
Set-TextValue "E12" "  -0.88%  "
Set-TextValue "D13" "1.641.81"
Set-TextValue "E13" "  -0.51%  "
Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.01"
Set-TextValue "E14" "  -1.64%  "
Set-TextValue "B15" "Polygon"
Set-TextValue "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.561"
Set-TextValue "E15" "  -1.05%  "
Set-TextValue "D16" "64.13"
Set-TextValue "E16" "  -2.48%  "
Set-TextValue "D17" "27.402.25"
Set-TextValue "E17" "  +0.20%  "
Set-TextValue "D18" "228.09"
Set-TextValue "E18" "  -2.60%  "
Set-TextValue "B19" "Chainlink"
Set-TextValue "C19" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "7.68"
Set-TextValue "E19" "  +3.20%  "
Set-TextValue "B20" "ShibaInu"
Set-TextValue "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D20" "0.0₃0722"
Set-TextValue "E20" "  -0.29%  "
Set-TextValue "E21" "  +0.20%  "
Set-TextValue "D22" "4.30"
Set-TextValue "E22" "  -2.17%  "
Set-TextValue "D23" "9.79"
Set-TextValue "E23" "  +5.55%  "
Set-TextValue "E24" "  -2.85%  "
Set-TextValue "D25" "149.80"
Set-TextValue "E25" "  +2.30%  "
Set-TextValue "D26" "6.94"
Set-TextValue "E26" "  -3.15%  "
Set-TextValue "E27" "  +1.67%  "
Set-TextValue "E28" "  +0.27%  "
Set-TextValue "D29" "15.52"
Set-TextValue "E29" "  -3.18%  "
Set-TextValue "E30" "  -0.63%  "
Set-TextValue "D31" "0.0487"
Set-TextValue "E31" "  -1.89%  "
Set-TextValue "E32" "  -0.56%  "
Set-TextValue "D33" "3.16"
Set-TextValue "E33" "  +2.13%  "
Set-TextValue "D34" "1.408.58"
Set-TextValue "E34" "  -3.60%  "
Set-TextValue "E35" "  +2.10%  "
Set-TextValue "E36" "  -2.07%  "
Set-TextValue "D37" "0.568"
Set-TextValue "E37" "  -0.88%  "
Set-TextValue "E38" "  -1.36%  "
Set-TextValue "D39" "0.871"
Set-TextValue "E39" "  -3.85%  "
Set-TextValue "D40" "0.900"
Set-TextValue "E40" "  +14.61%  "
Set-TextValue "E41" "  -0.41%  "
Set-TextValue "E42" "  +0.09%  "
Set-TextValue "B43" "mCoin"
Set-TextValue "C43" "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue "D43" "2.48"
Set-TextValue "E43" "  -0.27%  "
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "5.50"
Set-TextValue "E44" "  +1.44%  "
Set-TextValue "B45" "MXToken"
Set-TextValue "C45" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D45" "2.24"
Set-TextValue "E45" "  +1.27%  "
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "64.58"
Set-TextValue "E46" "  -0.70%  "
Set-TextValue "B47" "RocketPoolETH"
Set-TextValue "C47" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D47" "1.775.16"
Set-TextValue "E47" "  -0.88%  "
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.67"
Set-TextValue "E48" "  -3.00%  "
Set-TextValue "B49" "Quant"
Set-TextValue "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "85.72"
Set-TextValue "E49" "  -2.84%  "
Set-TextValue "B50" "BabyDogeCoin"
Set-TextValue "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D50" "0.0₆0106"
Set-TextValue "E50" "  +0.32%  "
Set-TextValue "B51" "Algorand"
Set-TextValue "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.0986"
Set-TextValue "E51" "  -2.02%  "
